$wb = $excel.ActiveWorkbook

# --- Sheet "Test Results": append rows 28-40 ---
$ws1 = $wb.Worksheets.Item("Test Results")

# Row 28
$ws1.Range("A28").Value = 'How many scriptures are there?'
$ws1.Range("B28").Value = 'basic-kn-index'
$ws1.Range("C28").Value = '13128_t091_mkmk_091'
$ws1.Range("D28").Value = 'यानि ये लोग वेद को कोट करते हैं कि इस वेद मंत्र के द्वारा हम ये सिद्धांत प्रस्तुत कर रहे हैं। पूर्व मीमांसा, इसके रचयिता जैमिनी हैं। ये 2 ग्रंथ इसमें बड़े बड़े हैं। 1 पूर्व मीमांसा और 1 कर्म मि मानसा। ये भारद्वाज का लिखा हुआ है। लेकिन कर्म मीमांसा का प्रचार अधिक नहीं है। जैमिनी के दर्शन का बहुत प्रचार है। ये क्या कहते हैं? ईश्वर वगैरह की आवश्यकता नहीं। ये ईश्वर, ईश्वर लोग बकते रहे। फिर कर्म करो क्या कर्म? यज्ञ किसके लिये? जब ईश्वर भी नहीं मानते। आप देवताओं के लिए। और देवता भी प्रधान नहीं हैं। मंत्र प्रधान है बेद। और उससे मिलेगा क्या स्वर्ग? तो उसमें तो दुख निवृति होगी। नहीं। अरे बड़ा'
$ws1.Range("E28").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13128_t091_mkmk_091.mp4'
$ws1.Range("F28").Value = '00:15:15,010'
$ws1.Range("G28").Value = '00:16:48,620'
$ws1.Range("H28").Value = 'N/A'
$ws1.Range("I28").Value = 'तीन प्रकार की वाणी होती है— विनिर्गत वाणी; स्मृत वाणी और कृत वाणी।'
$ws1.Range("J28").Value = 'N/A'
$ws1.Range("K28").Value = '00:39:00,100'
$ws1.Range("L28").Value = '00:39:09,600'
$ws1.Range("M28").Value = 'Fail'

# Row 29
$ws1.Range("A29").Value = 'कितने प्रकार के शास्त्र हैं?'
$ws1.Range("B29").Value = 'basic-kn-index'
$ws1.Range("C29").Value = '13092_t055_mkmk_055'
$ws1.Range("D29").Value = 'ने 100 अट्ठासी उपनिषदों को बताया है, लिखा है। और इनके नंबर में भी श्रीराम शास्त्री के कुछ भेद है। हमारे यहां विद्वानों, ने, आचार्यों ने अपने अपने हिसाब से नंबर किए हैं। इसलिए कुछ मत बहुमत वैश्म है। वेदांत में देख लो जिसको ब्रह्म सूत्र कहते हैं। शारीरिक भाष्य कहते हैं। तो सबसे इम्पोर्टेंट। हमारे यहाँ छे दर्शन होते हैं। आस्तिक दर्शन कहलाते हैं? वेदांत न्याय, सान, मिमांसा, पा, तंजलवैशेशिक, इनमें और दर्शन जो हैं वो तो तत्वज्ञान बताते हैं और दुःख निवृत्ति उनका उद्देश्य है। उनका विशेष महत्व नहीं है। केवल 1 वेदांत शास्त्र हमारे यहाँ ऐसा है जो आनंद प्राप्ति का लक्ष्य बताता है और'
$ws1.Range("E29").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13092_t055_mkmk_055.mp4'
$ws1.Range("F29").Value = '00:09:03,510'
$ws1.Range("G29").Value = '00:10:39,400'
$ws1.Range("H29").Value = 'N/A'
$ws1.Range("I29").Value = 'तीन प्रकार की वाणी होती है— विनिर्गत वाणी; स्मृत वाणी और कृत वाणी।'
$ws1.Range("J29").Value = 'N/A'
$ws1.Range("K29").Value = '00:39:00,100'
$ws1.Range("L29").Value = '00:39:09,600'
$ws1.Range("M29").Value = 'Fail'

# Row 30
$ws1.Range("A30").Value = 'कितने शरीर हैं हमारे पास?'
$ws1.Range("B30").Value = 'basic-kn-index'
$ws1.Range("C30").Value = '13051_t014_main_kaun_mera_kaun_014'
$ws1.Range("D30").Value = 'पहले कैसे लोग रिकॉर्ड करते थे, तो कहते थे महाराज जी, हमारा कैसे खत्म हो जायेगा? तो हम ऐसे हाथ उठा देंगे। तो आप तो थोड़ी देर रुक जाइएगा। यानी कुछ गोली ली निकालिए कुछ कीजिए। तब तक मैं कैचर बदल लूं। तो हमको वैसे मारना पड़ता था, नहीं तो वो उनके टेप में आगे नहीं। और ये फिलोसॉफी तो ऐसा सब्जेक्ट है। 1 शब्द भी गड़बड़ कर, 2 तो सब उल्टा हो जाएगा। अब देखो न। 1 शब्द की तो बात है। 1 कहता है? ब्रह्म निराकार ही होता है, बस लड़ाई होगी, ब्रह्म साकार ही होता है। अब असल में क्या है ही न, कहो भी कह। 2 झगड़ा ब्रह्म निराकार भी होता है, ब्रह्म साकार भी होता है। 1 शब्द का मामला सब उल्टा हो जाता है। तो'
$ws1.Range("E30").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13051_t014_main_kaun_mera_kaun_014.mp4'
$ws1.Range("F30").Value = '00:41:19,650'
$ws1.Range("G30").Value = '00:42:48,350'
$ws1.Range("H30").Value = 'N/A'
$ws1.Range("I30").Value = 'तीन शरीर होते हैं; हम लोगों के | इसका नाम—  स्थूल शरीर; और एक होता है सूक्ष्म शरीर; और एक होता है कारण शरीर।'
$ws1.Range("J30").Value = 'N/A'
$ws1.Range("K30").Value = '00:21:48,900'
$ws1.Range("L30").Value = '00:22:03,400'
$ws1.Range("M30").Value = 'Fail'

# Row 31
$ws1.Range("A31").Value = 'मैं शरीर हूँ?'
$ws1.Range("B31").Value = 'basic-kn-index'
$ws1.Range("C31").Value = '13105_t068_mkmk_068'
$ws1.Range("D31").Value = 'ये शरीर का। तुम शरीर तो हो नहीं, काय को फील करते हो। तुम्हारे शरीर को गाली दे रहा है। तुमको तो नहीं दे रहा है। तुम्हारे मकान की ओर। कोई चप्पल मार रहा है। तो क्यों रोते चिल्लाते हो। आज मैं अपने को देही मानता हूँ। इसलिए तो शिशुपाल ने जब सौ गाली दी, भगवान गिन रहे थे। 1, 2, 3, 76, 77, अठत्तर उन्यासी 80। जब सौ पूरे हो गए, तो चक्कर बुलाया और सिर काट दिया। अब सारे राजा राजाओं के सामने ये नाटक हुआ। राजाओं ने देखा, हाँ, ये तो भगवान है। और वो निकली आत्मा और श्री कृष्ण में मिल गई। युधिष्ठर वहीं बैठे थे, उन्होंने अकेले में भगवान से प्रश्न किया। महाराज वो आपमें मिल गया। उसको तो नरक'
$ws1.Range("E31").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13105_t068_mkmk_068.mp4'
$ws1.Range("F31").Value = '00:47:29,510'
$ws1.Range("G31").Value = '00:49:06,380'
$ws1.Range("H31").Value = 'N/A'
$ws1.Range("I31").Value = 'TBD'
$ws1.Range("J31").Value = 'N/A'
$ws1.Range("K31").Value = '00:31:23,700'
$ws1.Range("L31").Value = '00:31:25,600'
$ws1.Range("M31").Value = 'Fail'

# Row 32
$ws1.Range("A32").Value = 'आप कुछ पागल हो? क्या बकवास कर रहे हो?'
$ws1.Range("B32").Value = 'basic-kn-index'
$ws1.Range("C32").Value = '13108_t071_mkmk_071'
$ws1.Range("D32").Value = 'भगवान का बहुत बड़ा पागल खाना है, क्योंकि उनकी ये पागलों की संख्या भी अनंत है। 10, 20 करोड़ अरब नहीं, अनंत जीव माया बद्ध हैं, आनंद पाने के लिए व्याकुल हैं। और चुंकी बेड के अनुसार प्रत्येक जीव अनादि है, आज है, सनातन है, अतएव सदा से आनंद के लिए पागल है। और चूंकी कोई जीव 1 क्षण को भी अकरमा नहीं रह सकता, इसलिए प्रत्येक जीव अनाधिकार से प्रत्येक क्षण, तैलधारा व दबि्षिन्न केवल आनंद प्राप्ति के हेतु ही पागल है। चिंतन मनन प्लानिंग प्रैक्टिस में लगा है। अगर आप कहें हैं, मुझे जो पागल हैं, वो पागल हैं। हाँ, वो भी पागल है। भगवान ने हमको पागल डिक्लेयर किया।'
$ws1.Range("E32").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13108_t071_mkmk_071.mp4'
$ws1.Range("F32").Value = '00:08:54,560'
$ws1.Range("G32").Value = '00:10:29,370'
$ws1.Range("H32").Value = 'N/A'
$ws1.Range("I32").Value = 'N/A'
$ws1.Range("J32").Value = 'N/A'
$ws1.Range("K32").Value = 'N/A'
$ws1.Range("L32").Value = 'N/A'
$ws1.Range("M32").Value = 'Fail'

# Row 33
$ws1.Range("A33").Value = 'इन्द्रियाँ भी सब यहीं पड़ी।'
$ws1.Range("B33").Value = 'basic-kn-index'
$ws1.Range("C33").Value = '13041_t004_main_kaun_mera_kaun_004'
$ws1.Range("D33").Value = 'हैं। सजाते सम्बन्ध क्या 1 जाति के हैं। हम दोनों भगवान भी चेतन, हम भी चेतन, व हमारा अंश ही, हम, उसके अंश वो, हमारा माता, पिता, भ्राता सब कुछ और सख् संबन्ध हमारा सखा है, हमको सदा शक्ति देता है जीवित रहने का। प्रत्येक इंद्रियों में, मन में, बुद्धि में, त कर्म करने की शक्ति देता है। केनोपनिषद कहता है य बचा नब्बुडुतमजेन बाग बिदके तदेव ब्रह्म तम बिधि, यन मन सान मनुतेजेनाहुर मनोमतम तदेव ब्रह्म पम बिद्धि या चक्छुसानपशयति, जे न चकछुमपश्यतितदेव ब्रह्म या छरोतरेननचुणोती, ये न सत्र तदेव ब्रह्म या प्राण न प्राणि। ये न प्राण प्रणीयते। तदेव ब्रह्म बिद्धि यानि'
$ws1.Range("E33").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13041_t004_main_kaun_mera_kaun_004.mp4'
$ws1.Range("F33").Value = '00:17:57,450'
$ws1.Range("G33").Value = '00:19:40,470'
$ws1.Range("H33").Value = 'N/A'
$ws1.Range("I33").Value = 'यद्यपि कुछ आइडिया होता है कि जो ‘मेरा’ नहीं है; वो ‘मैं’ होगा।'
$ws1.Range("J33").Value = 'N/A'
$ws1.Range("K33").Value = '7:22.53'
$ws1.Range("L33").Value = '7:34.64'
$ws1.Range("M33").Value = 'Fail'

# Row 34
$ws1.Range("A34").Value = 'यद्यपि कुछ आइडिया होता है कि जो ‘मेरा’ नहीं है; वो ‘मैं’ होगा।'
$ws1.Range("B34").Value = 'basic-kn-index'
$ws1.Range("C34").Value = '13038_t001_main_kaun_mera_kaun_001'
$ws1.Range("D34").Value = 'लेकिन आत्मा को उसके लिए दुनिया में कोई खुर्दबीन बनी है। न बनेगी। न। तो दुनिया में कोई बुद्धिमान पैदा हुआ था। न होगा। जो मैं को जान लो। बुद्धि से आईडिया से आईडिया, धुआ आ गए। आइडिया कहते हैं इसको यत्र, यत्र धूमस तततत्रागनी। लेकिन यहाँ आइडिया काम नहीं देगा। न प्रत्यक्ष न। अनुमान। तो तीसरा प्रमाण है। शब्द प्रमाण। शब्द प्रमाण से जाना जायेगा। मैं कौन हूँ, मेरा कौन है शब्द प्रमाण? किसे कहते हैं? शब्द प्रमाण। माने। ऐसा प्रमाण जिसके शब्द पर कोई डाउट न हो। तो ऐसे 2 होते हैं। 1 तो जो मय को जान चूका है। ऐसा'
$ws1.Range("E34").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/10-2024/13038_t001_main_kaun_mera_kaun_001.mp4'
$ws1.Range("F34").Value = '00:32:21,750'
$ws1.Range("G34").Value = '00:33:54,740'
$ws1.Range("H34").Value = 'N/A'
$ws1.Range("I34").Value = 'यद्यपि कुछ आइडिया होता है कि जो ‘मेरा’ नहीं है; वो ‘मैं’ होगा।'
$ws1.Range("J34").Value = 'N/A'
$ws1.Range("K34").Value = '7:22.53'
$ws1.Range("L34").Value = '7:34.64'
$ws1.Range("M34").Value = 'Fail'

# Row 35
$ws1.Range("A35").Value = 'भगवान कहाँ हैं?'
$ws1.Range("B35").Value = 'basic-kn-index'
$ws1.Range("C35").Value = '13061_t024_main_kaun_mera_kaun_024'
$ws1.Range("D35").Value = 'उस किरण में कहाँ से आई? सूर्य से? तो फिर हमारी हैसियत, क्या है? हमारी तो पर्सनैलिटी भगवान पर डिपेंड करती है। इसलिए वेद कहता है सकारण करना। दीपा, दीपा, स्वेता, शतरोपनिशचछानौ प्रधान के तरगपतिरगुणशा श्वेता चतरो परिषद छे 16 प्रधान पुरुष स्वरा भगवान, हे जीव और माया का अध्यक्ष है, गवर्नर है। उसके बिना न जीव का जीवत्व है। न माया बिचारी जड़ वस्त का जड़त्व हैं, लेकिन उसकी शक्ति पाकर देखो, माया का हाल सो माया सब जगहीं नचावा शिव चतुरानन देखी, डेरा ही अपर जीव के ही लेके माही शिव बिरनचिकामोहई को है'
$ws1.Range("E35").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13061_t024_main_kaun_mera_kaun_024.mp4'
$ws1.Range("F35").Value = '00:29:15,380'
$ws1.Range("G35").Value = '00:30:37,300'
$ws1.Range("H35").Value = 'N/A'
$ws1.Range("I35").Value = 'भगवान् ने कहा—-“हाँ; रे! मैं सब जगह रहता हूँ; देख! देख ले आँख से; तेरा काल हूँ। मैं...मैं राक्षस-पाक्षस नहीं मानता; गन्दा-फंदा नहीं मानता; मैं सर्वत्र रहता हूँ।”'
$ws1.Range("J35").Value = 'N/A'
$ws1.Range("K35").Value = '0:35:30'
$ws1.Range("L35").Value = '0:35:53'
$ws1.Range("M35").Value = 'Fail'

# Row 36
$ws1.Range("A36").Value = 'कितने प्रकार के वाणी होते हैं?'
$ws1.Range("B36").Value = 'basic-kn-index'
$ws1.Range("C36").Value = '13061_t024_main_kaun_mera_kaun_024'
$ws1.Range("D36").Value = 'कोई व्यक्ति कहीं भी, कभी भी रा शब्द को बोलता है। raबडाrातो मेरा रोम रोम रोमांचित हो जाता है। हर्ष के साथ कभी यह माँ बोलेगा, भले ही उसकी मंशा रावण। बोलने के यही भगवन नाम हैं, जो हम लोग संकीर्तन करते हैं। इसी के विषय में भगवान शंकर कह रहे हैं। अतो भगवन नाम, संकीर्तन आदि में बहुत अधिक मन का संयोग करके उसकी इंपार्टेंस को रियलाइज करते हुए बोलना चाहिए। भोग गोबिंद गो कला गोला गो गो पाला a बोली वृन्दावन बिहारी लाल की अब आप लोग सावधान हो जायें। मैं कौन मेरा कौन इन प्रश्नों के समाधान में अब तक आप'
$ws1.Range("E36").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13061_t024_main_kaun_mera_kaun_024.mp4'
$ws1.Range("F36").Value = '00:06:54,180'
$ws1.Range("G36").Value = '00:11:07,440'
$ws1.Range("H36").Value = 'N/A'
$ws1.Range("I36").Value = 'तीन प्रकार की वाणी होती है— विनिर्गत वाणी; स्मृत वाणी और कृत वाणी।'
$ws1.Range("J36").Value = 'N/A'
$ws1.Range("K36").Value = '39:00.1'
$ws1.Range("L36").Value = '39:09.6'
$ws1.Range("M36").Value = 'Fail'

# Row 37
$ws1.Range("A37").Value = 'मैं जीव को देख सकता हूँ?'
$ws1.Range("B37").Value = 'basic-kn-index'
$ws1.Range("C37").Value = '13108_t071_mkmk_071'
$ws1.Range("D37").Value = 'पड़ते हैं यहाँ आप लोग बैठे हैं कितने उड़ रहे हैं जीव कहाँ उड़ रहे हैं अरे ये आँख नहीं देख रही हैं, खुर्दबीन लगाकर देखो तो मालूम पड़ जाए वो जा रहे हैं पेट में आपके आप बोलते हैं, खाते पीते हैं सब अरे पेट में आपके तमाम जीव हैं, खून देखा है हाँ उसमें जीव दिखते हैं जरा खुद बी सी देखो हर जगह जी तो संसार के पंच महाभूत के समान भी सब नहीं देख सकते हैं आँख से तो भगवान की बात क्या करते हैं भगवान को तो वो आपके बेटे से भी खराब दिखाई पड़ेगा तो आप कहेंगे यही भगवान है धन्य हैं तुम्हारे भगवान, गुरु जी तुमको भी धन्य हैं दोनों को नमस्कार कल ऐसी नहीं आऊंगा मैं तो समझ रहा था की अनन्त कोट कन्दर्प दर्प दलन पतियान हैं,'
$ws1.Range("E37").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13108_t071_mkmk_071.mp4'
$ws1.Range("F37").Value = '00:59:35,560'
$ws1.Range("G37").Value = '01:00:58,560'
$ws1.Range("H37").Value = 'N/A'
$ws1.Range("I37").Value = 'इतना सूक्ष्म है ‘मैं’ कि किसी यंत्र से कोई नहीं देख सका; न देख सकेगा'
$ws1.Range("J37").Value = 'N/A'
$ws1.Range("K37").Value = '31:23.7'
$ws1.Range("L37").Value = '31:25.6'
$ws1.Range("M37").Value = 'Fail'

# Row 38
$ws1.Range("A38").Value = 'सारे जीवों का क्या प्रश्न है?'
$ws1.Range("B38").Value = 'basic-kn-index'
$ws1.Range("C38").Value = '13038_t001_main_kaun_mera_kaun_001'
$ws1.Range("D38").Value = 'हैं। वो प्रश्न कौन हैं? हम कौन हैं? हमारा कौन है? बस, हम कौन हैं? आपको ये प्रश्न सुनकर कुछ आश्यर्य अवश्य होगा। इसका उत्तर तो सब जानते हैं। नट। कोई कोई जानता है, करोड़ो में कोई जानता है। हम लोगो को भी जानना है। वेद कहता है हे द बेदी, दत, सत्य मस्त नच्चे दिहा बेदिन महति विनष्टि केनो पनिषद 25 अरे मनुष्यों। इस प्रश्न का उत्तर समझना बहुत आवश्यक है। अगर नहीं समझोगे तो बहुत बड़ी हानी होगी। क्या हानी होगी?'
$ws1.Range("E38").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/10-2024/13038_t001_main_kaun_mera_kaun_001.mp4'
$ws1.Range("F38").Value = '00:07:23,710'
$ws1.Range("G38").Value = '00:09:21,490'
$ws1.Range("H38").Value = 'N/A'
$ws1.Range("I38").Value = 'समस्त विश्व के समस्त प्राणी केवल एक प्रश्न के हल करने में लगे हैं; वो  प्रश्न कौन है? हम कौन हैं? हमारा कौन है?'
$ws1.Range("J38").Value = 'N/A'
$ws1.Range("K38").Value = '00:07:09,500'
$ws1.Range("L38").Value = '00:07:38,800'
$ws1.Range("M38").Value = 'Fail'

# Row 39
$ws1.Range("A39").Value = 'हमारे पास कितने विभिन्न प्रकार के शरीर होते हैं, और वे किस प्रकार से एक-दूसरे से भिन्न होते हैं?'
$ws1.Range("B39").Value = 'basic-kn-index'
$ws1.Range("C39").Value = '13136_t099_mkmk_099'
$ws1.Range("D39").Value = 'में भेद भी होता है, अभेद भी होता है, अभेद यों है को, वह भी चेतन हैं, हम भी चेतन है, वह भी ज्ञान युक्त हैं, हम भी ज्ञान युक्त हैं, वह भी आनन्द सिंधु हैं और हम भी उसी आनन्द सिंधु के इच्छुक हैं। अर्थात हम कुछ अंश में उससे अभिन्न हैं, अधिक अंशों में भिन्न हैं। अनेक वेद मंत्रों के द्वारा आपको बताया गया है कि वो ब्रह्म इस जीव माया से पृथक हैं। वो सर्व शक्तिमान हैं, ये अल्प शक्तिमान हैं, वो सर्वग्य है, यह अल्पज्ञ है, वो मायाधीश हैं, ये मायाधीन हैं, अनेक अंतर है। वेदांत में बड़ा विस्तृत विचार किया गया है। यथा'
$ws1.Range("E39").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13136_t099_mkmk_099.mp4'
$ws1.Range("F39").Value = '00:12:32,740'
$ws1.Range("G39").Value = '00:14:07,480'
$ws1.Range("H39").Value = 'N/A'
$ws1.Range("I39").Value = 'तीन शरीर होते हैं; हम लोगों के | इसका नाम—  स्थूल शरीर; और एक होता है सूक्ष्म शरीर; और एक होता है कारण शरीर।'
$ws1.Range("J39").Value = 'N/A'
$ws1.Range("K39").Value = '00:21:48,900'
$ws1.Range("L39").Value = '00:22:03,400'
$ws1.Range("M39").Value = 'Fail'

# Row 40
$ws1.Range("A40").Value = 'आत्मा के शरीर में बाहर से प्रवेश करने, आने और जाने को समझाने के लिए महाराज जी कौन सा उदाहरण देते हैं'
$ws1.Range("B40").Value = 'basic-kn-index'
$ws1.Range("C40").Value = '13137_t100_mkmk_100'
$ws1.Range("D40").Value = 'है। ये शरीर चलाने के लिए भगवान ने बनाया है। क्योंकि शरीर के बिना साधना नहीं हो सकती। तनु बिनु, भजन बेद नहीं बरना। तो ये अवतार का कारण नणामनश्रेजसारथाय। व्यक्त भगवतो निरगुण गुणात्म। 10 वनते, 14 विद्या ने फिर कहा, तेजी पर कल्याण करने के लिए भगवान जन्म लेते हैं और कर्म करते हैं। इसलिए गीता ने कहा जन्म कर्म चमे दिव्यम। अर्जुन, मेरा जन्म भी दिव्य है, कर्म भी दिव्य है और तेरा कर्म दिव्य है। जन्म डिब्बेमे देखो। 3 प्रकार का देह होता है देह। ये शरीर 1 पंच महाभूत का और 1 चिदानंद में देह दिव्य देह और 1 देही'
$ws1.Range("E40").Value = 'https://jkp-llm.s3.ap-south-1.amazonaws.com/Videos/11-2024/13137_t100_mkmk_100.mp4'
$ws1.Range("F40").Value = '00:51:39,620'
$ws1.Range("G40").Value = '00:53:28,990'
$ws1.Range("H40").Value = 'N/A'
$ws1.Range("I40").Value = 'एक कमरा है; उसमें अँधेरा है। हमने बाहर से दीपक जलाकरके कमरे में रख दिया।तो; दीपक स्वयं प्रकाश है और उसने कमरे में भी प्रकाश कर दिया। दीपक को हमने दूसरे कमरे में लेजाके रख दिया। “तो?” तो कमरे में अँधेरा हो गया। ऐसे ही; ये जीव है; बाहर से इस शरीर में माँ के पेट में आया। आया; और फिर चला गया'
$ws1.Range("J40").Value = 'N/A'
$ws1.Range("K40").Value = '00:19:21.86'
$ws1.Range("L40").Value = '00:20:16.86'
$ws1.Range("M40").Value = 'Fail'

# --- Sheet "Metrics": append rows 28-40 ---
$ws2 = $wb.Worksheets.Item("Metrics")

# Row 28
$ws2.Range("A28").Value = 'basic-kn-index'
$ws2.Range("B28").Value = 'similarity-search'
$ws2.Range("C28").Value = 0

# Row 29
$ws2.Range("A29").Value = 'basic-kn-index'
$ws2.Range("B29").Value = 'similarity-search'
$ws2.Range("C29").Value = 0

# Row 30
$ws2.Range("A30").Value = 'basic-kn-index'
$ws2.Range("B30").Value = 'similarity-search'
$ws2.Range("C30").Value = 0

# Row 31
$ws2.Range("A31").Value = 'basic-kn-index'
$ws2.Range("B31").Value = 'similarity-search'
$ws2.Range("C31").Value = 0

# Row 32
$ws2.Range("A32").Value = 'basic-kn-index'
$ws2.Range("B32").Value = 'similarity-search'
$ws2.Range("C32").Value = 0

# Row 33
$ws2.Range("A33").Value = 'basic-kn-index'
$ws2.Range("B33").Value = 'similarity-search'
$ws2.Range("C33").Value = 0

# Row 34
$ws2.Range("A34").Value = 'basic-kn-index'
$ws2.Range("B34").Value = 'similarity-search'
$ws2.Range("C34").Value = 0

# Row 35
$ws2.Range("A35").Value = 'basic-kn-index'
$ws2.Range("B35").Value = 'similarity-search'
$ws2.Range("C35").Value = 0

# Row 36
$ws2.Range("A36").Value = 'basic-kn-index'
$ws2.Range("B36").Value = 'similarity-search'
$ws2.Range("C36").Value = 0

# Row 37
$ws2.Range("A37").Value = 'basic-kn-index'
$ws2.Range("B37").Value = 'similarity-search'
$ws2.Range("C37").Value = 0

# Row 38
$ws2.Range("A38").Value = 'basic-kn-index'
$ws2.Range("B38").Value = 'similarity-search'
$ws2.Range("C38").Value = 0

# Row 39
$ws2.Range("A39").Value = 'basic-kn-index'
$ws2.Range("B39").Value = 'similarity-search'
$ws2.Range("C39").Value = 0

# Row 40
$ws2.Range("A40").Value = 'basic-kn-index'
$ws2.Range("B40").Value = 'similarity-search'
$ws2.Range("C40").Value = 0
